# Updates the crypto price/volume table to the latest scraped values.
# D-column "price" cells must stay plain text (they use dotted thousand
# separators like "61.847.52" and some would otherwise be auto-parsed as
# numbers by Excel, e.g. "72.68" or "1.00"). To keep them as text while
# leaving the cell style untouched (no explicit style index in the
# original file), we temporarily force a text number format, set the
# value, and then restore the "Normal" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "61.847.52"
$ws.Range("E2").Value = "  +1.27%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.460.04"
$ws.Range("E3").Value = "  +2.74%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
Set-TextValue "D5" "581.02"
$ws.Range("E5").Value = "  +1.74%  "

# Row 6 - Solana
Set-TextValue "D6" "148.98"
$ws.Range("E6").Value = "  +9.71%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.460.99"
$ws.Range("E7").Value = "  +2.81%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - XRP
Set-TextValue "D9" "0.475"
$ws.Range("E9").Value = "  +1.27%  "

# Row 10 - Toncoin
Set-TextValue "D10" "7.78"
$ws.Range("E10").Value = "  +3.77%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.74%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.392"
$ws.Range("E12").Value = "  +1.99%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "4.051.26"
$ws.Range("E13").Value = "  +2.79%  "

# Row 14 - Avalanche
Set-TextValue "D14" "28.29"
$ws.Range("E14").Value = "  +8.79%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  -0.48%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +1.88%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "3.456.86"
$ws.Range("E17").Value = "  +2.71%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "61.868.61"
$ws.Range("E18").Value = "  +1.11%  "

# Row 19 - Polkadot
Set-TextValue "D19" "6.36"
$ws.Range("E19").Value = "  +9.54%  "

# Row 20 - Chainlink
Set-TextValue "D20" "14.38"
$ws.Range("E20").Value = "  +2.82%  "

# Row 21 - Uniswap
Set-TextValue "D21" "9.47"
$ws.Range("E21").Value = "  +2.41%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "385.74"
$ws.Range("E22").Value = "  +2.48%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.571"
$ws.Range("E23").Value = "  +3.12%  "

# Row 24 - WrappedeETH
Set-TextValue "D24" "3.595.82"
$ws.Range("E24").Value = "  +2.85%  "

# Row 25 - was Dai, now Litecoin (rows 25 and 27 swap coins)
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D25" "72.68"
$ws.Range("E25").Value = "  +2.35%  "

# Row 26 - LEO
$ws.Range("E26").Value = "  +0.98%  "

# Row 27 - was Litecoin, now Dai
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D27" "1.00"
$ws.Range("E27").Value = "  -0.21%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  -1.44%  "

# Row 29 - Kaspa
Set-TextValue "D29" "0.181"
$ws.Range("E29").Value = "  +9.28%  "

# Row 30 - RenderToken
Set-TextValue "D30" "7.85"
$ws.Range("E30").Value = "  +4.68%  "

# Row 31 - Binance-PegBSC-USD
Set-TextValue "D31" "1.00"
$ws.Range("E31").Value = "  -0.52%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  -13.20%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "8.26"
$ws.Range("E33").Value = "  +1.40%  "

# Row 34 - PancakeSwap
Set-TextValue "D34" "2.18"
$ws.Range("E34").Value = "  +1.85%  "

# Row 35 - USDe (no change)

# Row 36 - EthereumClassic
Set-TextValue "D36" "23.97"
$ws.Range("E36").Value = "  +1.83%  "

# Row 37 - Aptos
$ws.Range("E37").Value = "  +4.48%  "

# Row 38 - NEARProtocol
Set-TextValue "D38" "5.23"
$ws.Range("E38").Value = "  +0.59%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +2.82%  "

# Row 40 - Monero
Set-TextValue "D40" "166.30"
$ws.Range("E40").Value = "  +0.77%  "

# Row 41 - Hedera
Set-TextValue "D41" "0.0792"
$ws.Range("E41").Value = "  +5.17%  "

# Row 42 - was Mantle, now EnergySwap (rows 42 and 43 swap coins)
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D42" "25.91"
$ws.Range("E42").Value = "  +8.38%  "

# Row 43 - was EnergySwap, now Mantle
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D43" "0.797"
$ws.Range("E43").Value = "  +3.68%  "

# Row 44 - was Stacks, now FirstDigitalUSD (rows 44 and 45 swap coins)
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D44" "1.00"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45 - was FirstDigitalUSD, now Stacks
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D45" "1.73"
$ws.Range("E45").Value = "  +1.65%  "

# Row 46 - OKB
Set-TextValue "D46" "42.37"
$ws.Range("E46").Value = "  +2.20%  "

# Row 47 - Filecoin
Set-TextValue "D47" "4.49"
$ws.Range("E47").Value = "  +2.39%  "

# Row 48 - ONDO
Set-TextValue "D48" "1.18"
$ws.Range("E48").Value = "  -1.91%  "

# Row 49 - Maker
Set-TextValue "D49" "2.605.41"
$ws.Range("E49").Value = "  +10.94%  "

# Row 50 - Cosmos
Set-TextValue "D50" "6.98"
$ws.Range("E50").Value = "  +2.56%  "

# Row 51 - InjectiveProtocol
Set-TextValue "D51" "23.37"
$ws.Range("E51").Value = "  +0.57%  "
